$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Amapá
$ws.Range("B2").Value = "Diferença 2024/03 - 2023/03"
$ws.Range("C2").Value = 4.425906615687637

# Row 3 - Pernambuco -> Bahia
$ws.Range("A3").Value = "Bahia"
$ws.Range("B3").Value = "Diferença 2024/03 - 2023/03"
$ws.Range("C3").Value = 3.623377883072379

# Row 4 - Bahia -> Pernambuco
$ws.Range("A4").Value = "Pernambuco"
$ws.Range("B4").Value = "Diferença 2024/03 - 2023/03"
$ws.Range("C4").Value = 2.696239622441539

# Row 5 - Piauí -> Ceará
$ws.Range("A5").Value = "Ceará"
$ws.Range("B5").Value = "Diferença 2024/03 - 2023/03"
$ws.Range("C5").Value = 2.503832437878714

# Row 6 - Tocantins -> Rio de Janeiro
$ws.Range("A6").Value = "Rio de Janeiro"
$ws.Range("B6").Value = "Diferença 2024/03 - 2023/03"
$ws.Range("C6").Value = 2.358712302206129

# Row 7 - Nordeste
$ws.Range("B7").Value = "Diferença 2024/03 - 2023/03"
$ws.Range("C7").Value = 2.178119867450107

# Row 8 - Sergipe
$ws.Range("B8").Value = "Diferença 2024/03 - 2023/03"
$ws.Range("C8").Value = 1.411302622922179
$ws.Range("D8").Value = "10º"

# Row 9 - Nordeste
$ws.Range("B9").Value = "Diferença 2024/03 - 2023/03"
$ws.Range("C9").Value = 2.178119867450107

# Row 10 - Brasil
$ws.Range("B10").Value = "Diferença 2024/03 - 2023/03"
$ws.Range("C10").Value = 1.326225861527845
